$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update query text (column B) for rows whose SLR query was replaced by SQW query
$ws.Range("B31").Value = """software development methodologies"" OR ""software project management"" OR (software process line) OR ""code review processes"" OR ""risk management in software"" OR ""software testing strategies"" OR ""software documentation"" OR ""agile software development"" OR ""process maturity model"" OR ""software quality assurance"" OR ""user acceptance testing"" OR ""agile methodologies"" OR ""software development lifecycle"" OR ""software product line engineering"" OR ""software engineering best practices"" OR ""devops practices"" OR ((""process modeling"" OR ""system design"" OR ""change management"" OR ""process optimization"" OR ""process automation"" OR ""maintenance processes"" OR ""performance metrics"" OR ""iterative development"" OR ""software architecture"" OR ""continuous integration"" OR ""requirements engineering"" OR ""software process improvement"" OR ""release management"" OR ""configuration management"" OR ""version control systems"" OR ""stakeholder engagement"") AND (Software Process))"
$ws.Range("B33").Value = """field-programmable gate arrays"" OR ""low-latency processing"" OR ""real-time analytics"" OR ""concept drift"" OR ""stream benchmarks"" OR ""workload variations"" OR ""reactive programming"" OR ""memory reuse"" OR ""floating-point units"" OR ""edge stream processing"" OR ""query latency"" OR (data stream processing latency) OR ""dataflow programming"" OR ""push-based streaming"" OR ""microbenchmarks"" OR ""rpc mechanisms"" OR ""dma management"" OR ""pull-based streaming"" OR ""scoreboarding"" OR ""ingestion/storage integration"" OR ((""sensor data"" OR ""distributed data processing"" OR ""parallel processing"" OR ""stream processing"" OR ""scalability"" OR ""adaptability"" OR ""kernel functions"" OR ""edge computing"" OR ""shared memory"" OR ""streaming algorithms"" OR ""bandwidth optimization"" OR ""compiler tools"" OR ""hardware acceleration"" OR ""graphics processing units"" OR ""pipelining"") AND (Processing Latency))"
$ws.Range("B35").Value = """business rules"" OR ""business process frameworks"" OR ""workflow management"" OR ""enterprise architecture"" OR ""lean management"" OR ""requirements gathering"" OR ""six sigma"" OR ""business process reengineering"" OR ""bpmn"" OR ""organizational process assets"" OR (business process meta models) OR ((""process analysis"" OR ""process improvement"" OR ""process simulation"" OR ""change management"" OR ""risk management"" OR ""knowledge management"" OR ""process optimization"" OR ""continuous improvement"" OR ""modeling languages"" OR ""process automation"" OR ""process mapping"" OR ""performance metrics"" OR ""business process management"" OR ""quality assurance"" OR ""stakeholder analysis"" OR ""business process modeling"" OR ""value chain analysis"" OR ""process governance"" OR ""process documentation"") AND (business metamodels))"
$ws.Range("B37").Value = """data locality"" OR ""system bottlenecks"" OR ""computational throughput"" OR ""memory bandwidth"" OR ""workload distribution"" OR ""core utilization"" OR ""simd"" OR ""interconnect bandwidth"" OR ""thread management"" OR ""multicore processors"" OR ""hardware threading"" OR ""multicore architecture"" OR ""task granularity"" OR ""resource contention"" OR ""cache coherence"" OR ""reuse profiles"" OR (multicore performance prediction) OR ""amdahl's law"" OR ((""performance modeling"" OR ""energy efficiency"" OR ""analytical modeling"" OR ""performance prediction"" OR ""latency"" OR ""performance metrics"" OR ""scheduling algorithms"" OR ""scalability"" OR ""load balancing"" OR ""parallel computing"" OR ""software optimization"" OR ""benchmarking"" OR ""task scheduling"") AND (Multicore Performance))"
$ws.Range("B39").Value = """virtualization"" OR ""cloud strategy"" OR ""cloud optimization"" OR ""disaster recovery"" OR ""hybrid cloud"" OR ""cloud scalability"" OR ""cloud transformation"" OR (cloud migration) OR ""cloud adoption"" OR ""cloud cost optimization"" OR ""cloud compliance"" OR ""cloud backup"" OR ""cloud native applications"" OR ""serverless computing"" OR ""devops in cloud"" OR ""iaas migration"" OR ""paas migration"" OR ""saas migration"" OR ((""data migration"" OR ""application migration"" OR ""cloud services"" OR ""cloud networking"" OR ""cloud performance"" OR ""cloud management"" OR ""cloud security"" OR ""cloud architecture"" OR ""multi-cloud"" OR ""cloud storage"" OR ""cloud monitoring"" OR ""cloud infrastructure"") AND (Cloud Migration))"
$ws.Range("B41").Value = """defect identification"" OR ""software quality"" OR ""software complexity"" OR ""change-proneness"" OR ""metric validation"" OR ""cost evaluation framework"" OR ""fault prediction models"" OR ""testing resource allocation"" OR ""fault-proneness"" OR ""source code metrics"" OR ""software maintainability"" OR ""software metrics suite"" OR ""internal software metrics"" OR (software fault prediction metrics) OR ""fault-prone classes"" OR ""apache ecosystem"" OR ""eclipse ecosystem"" OR ""sonarqube rules"" OR ((""predictive models"" OR ""performance evaluation"" OR ""statistical methods"" OR ""data quality"" OR ""reliability"" OR ""empirical study"" OR ""comments"" OR ""statistical correlation"" OR ""machine learning"" OR ""classification techniques"" OR ""prediction accuracy"" OR ""inheritance"" OR ""open-source systems"" OR ""ensemble methods"" OR ""cohesion"" OR ""ensemble learning"" OR ""generalizability"") AND (software error prediction))"
$ws.Range("B43").Value = """software maintenance"" OR ""software metrics"" OR ""model transparency"" OR ""ensemble classification"" OR ""defect localization"" OR ""effort metrics"" OR ""software artifacts"" OR ""defect prediction models"" OR ""quality assurance resources"" OR ""heterogeneous classifiers"" OR (software defect prediction) OR ""software sizing metrics"" OR ""model explainability"" OR ""post-release defects"" OR ""static code metrics"" OR ""historical defect information"" OR ""cyclomatic complexity"" OR ""defect-prone modules"" OR ""source code similarity metrics"" OR ""cross-company project metrics"" OR ((""predictive models"" OR ""cost model"" OR ""testing costs"" OR ""machine learning"" OR ""defectiveness"" OR ""software quality"" OR ""resource allocation"" OR ""defect density"" OR ""quality assurance"" OR ""defect repair"" OR ""project-level information"") AND (Software defect))"

# Update numeric metric columns (H, I, J, K, and for rows 30/32/34/36/40/42 also C-G)
$ws.Cells.Item(2, 8).Value = 0.128
$ws.Cells.Item(2, 9).Value = 0.38
$ws.Cells.Item(2, 10).Value = 0.122
$ws.Cells.Item(2, 11).Value = 0.37
$ws.Cells.Item(3, 8).Value = 0.835
$ws.Cells.Item(3, 10).Value = 0.801
$ws.Cells.Item(3, 11).Value = 0.9399999999999999
$ws.Cells.Item(4, 8).Value = 0.362
$ws.Cells.Item(4, 9).Value = 0.16
$ws.Cells.Item(4, 10).Value = 0.216
$ws.Cells.Item(5, 8).Value = 0.654
$ws.Cells.Item(5, 9).Value = 0.5600000000000001
$ws.Cells.Item(5, 10).Value = 0.368
$ws.Cells.Item(6, 8).Value = 0.477
$ws.Cells.Item(6, 9).Value = 0.06
$ws.Cells.Item(6, 10).Value = 0.295
$ws.Cells.Item(6, 11).Value = 0.23
$ws.Cells.Item(7, 8).Value = 0.738
$ws.Cells.Item(7, 9).Value = 0.8
$ws.Cells.Item(7, 10).Value = 0.529
$ws.Cells.Item(7, 11).Value = 0.76
$ws.Cells.Item(8, 8).Value = 0.822
$ws.Cells.Item(8, 10).Value = 0.62
$ws.Cells.Item(8, 11).Value = 0.01
$ws.Cells.Item(9, 8).Value = 0.185
$ws.Cells.Item(9, 9).Value = 0.42
$ws.Cells.Item(9, 10).Value = 0.11
$ws.Cells.Item(9, 11).Value = 0.34
$ws.Cells.Item(10, 8).Value = 0.15
$ws.Cells.Item(10, 9).Value = 0.33
$ws.Cells.Item(10, 10).Value = 0.07099999999999999
$ws.Cells.Item(10, 11).Value = 0.24
$ws.Cells.Item(12, 8).Value = 0.573
$ws.Cells.Item(12, 10).Value = 0.476
$ws.Cells.Item(12, 11).Value = 0.13
$ws.Cells.Item(13, 8).Value = 0.634
$ws.Cells.Item(13, 10).Value = 0.493
$ws.Cells.Item(13, 11).Value = 0.72
$ws.Cells.Item(14, 8).Value = 0.06900000000000001
$ws.Cells.Item(14, 9).Value = 0.24
$ws.Cells.Item(14, 10).Value = 0.028
$ws.Cells.Item(14, 11).Value = 0.12
$ws.Cells.Item(18, 8).Value = 0.513
$ws.Cells.Item(18, 9).Value = 0.15
$ws.Cells.Item(18, 10).Value = 0.343
$ws.Cells.Item(18, 11).Value = 0.32
$ws.Cells.Item(20, 8).Value = 0.977
$ws.Cells.Item(20, 10).Value = 0.88
$ws.Cells.Item(21, 8).Value = 0.926
$ws.Cells.Item(21, 9).Value = 0.47
$ws.Cells.Item(21, 10).Value = 0.798
$ws.Cells.Item(21, 11).Value = 0.54
$ws.Cells.Item(24, 8).Value = 0.365
$ws.Cells.Item(24, 9).Value = 0.19
$ws.Cells.Item(24, 10).Value = 0.186
$ws.Cells.Item(25, 8).Value = 0.342
$ws.Cells.Item(25, 9).Value = 0.72
$ws.Cells.Item(25, 10).Value = 0.182
$ws.Cells.Item(25, 11).Value = 0.53
$ws.Cells.Item(26, 8).Value = 0.264
$ws.Cells.Item(26, 9).Value = 0.26
$ws.Cells.Item(26, 10).Value = 0.162
$ws.Cells.Item(27, 8).Value = 0.082
$ws.Cells.Item(27, 9).Value = 0.31
$ws.Cells.Item(27, 10).Value = 0.024
$ws.Cells.Item(27, 11).Value = 0.11
$ws.Cells.Item(28, 8).Value = 0.434
$ws.Cells.Item(28, 9).Value = 0.15
$ws.Cells.Item(28, 10).Value = 0.332
$ws.Cells.Item(28, 11).Value = 0.25
$ws.Cells.Item(29, 8).Value = 0.58
$ws.Cells.Item(29, 10).Value = 0.441
$ws.Cells.Item(30, 3).Value = 0.605
$ws.Cells.Item(30, 4).Value = 0.02
$ws.Cells.Item(30, 5).Value = 0.09
$ws.Cells.Item(30, 6).Value = 0.316
$ws.Cells.Item(30, 7).Value = 0.22
$ws.Cells.Item(30, 8).Value = 0.701
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0.261
$ws.Cells.Item(30, 11).Value = 0.3
$ws.Cells.Item(31, 8).Value = 0.749
$ws.Cells.Item(31, 9).Value = 0.66
$ws.Cells.Item(31, 10).Value = 0.298
$ws.Cells.Item(31, 11).Value = 0.6
$ws.Cells.Item(32, 3).Value = 0.261
$ws.Cells.Item(32, 4).Value = 0.016
$ws.Cells.Item(32, 5).Value = 0.06
$ws.Cells.Item(32, 6).Value = 0.014
$ws.Cells.Item(32, 7).Value = 0.06
$ws.Cells.Item(32, 8).Value = 0.034
$ws.Cells.Item(32, 9).Value = 0.14
$ws.Cells.Item(32, 10).Value = 0.026
$ws.Cells.Item(32, 11).Value = 0.11
$ws.Cells.Item(33, 8).Value = 0.271
$ws.Cells.Item(33, 9).Value = 0.65
$ws.Cells.Item(33, 10).Value = 0.222
$ws.Cells.Item(33, 11).Value = 0.59
$ws.Cells.Item(34, 3).Value = 0.615
$ws.Cells.Item(34, 4).Value = 0.061
$ws.Cells.Item(34, 5).Value = 0.21
$ws.Cells.Item(34, 6).Value = 0.107
$ws.Cells.Item(34, 7).Value = 0.26
$ws.Cells.Item(34, 8).Value = 0.152
$ws.Cells.Item(34, 9).Value = 0.39
$ws.Cells.Item(34, 10).Value = 0.112
$ws.Cells.Item(34, 11).Value = 0.34
$ws.Cells.Item(35, 8).Value = 0.38
$ws.Cells.Item(35, 9).Value = 0.75
$ws.Cells.Item(35, 10).Value = 0.26
$ws.Cells.Item(35, 11).Value = 0.64
$ws.Cells.Item(36, 3).Value = 0.273
$ws.Cells.Item(36, 4).Value = 0.118
$ws.Cells.Item(36, 5).Value = 0.2
$ws.Cells.Item(36, 7).Value = 0.01
$ws.Cells.Item(36, 8).Value = 0.158
$ws.Cells.Item(36, 9).Value = 0.4
$ws.Cells.Item(36, 10).Value = 0.1
$ws.Cells.Item(36, 11).Value = 0.32
$ws.Cells.Item(40, 3).Value = 0.5620000000000001
$ws.Cells.Item(40, 4).Value = 0.078
$ws.Cells.Item(40, 5).Value = 0.24
$ws.Cells.Item(40, 6).Value = 1
$ws.Cells.Item(40, 7).Value = 0.23
$ws.Cells.Item(40, 8).Value = 0.254
$ws.Cells.Item(40, 9).Value = 0.53
$ws.Cells.Item(40, 10).Value = 0.15
$ws.Cells.Item(40, 11).Value = 0.43
$ws.Cells.Item(41, 8).Value = 0.248
$ws.Cells.Item(41, 9).Value = 0.62
$ws.Cells.Item(41, 10).Value = 0.138
$ws.Cells.Item(41, 11).Value = 0.44
$ws.Cells.Item(42, 3).Value = 0.519
$ws.Cells.Item(42, 4).Value = 0.142
$ws.Cells.Item(42, 5).Value = 0.32
$ws.Cells.Item(42, 6).Value = 0.155
$ws.Cells.Item(42, 7).Value = 0.31
$ws.Cells.Item(42, 8).Value = 0.603
$ws.Cells.Item(42, 9).Value = 0.51
$ws.Cells.Item(42, 10).Value = 0.402
$ws.Cells.Item(42, 11).Value = 0.57
$ws.Cells.Item(43, 8).Value = 0.724
$ws.Cells.Item(43, 9).Value = 0.9
$ws.Cells.Item(43, 10).Value = 0.635
$ws.Cells.Item(43, 11).Value = 0.87
